# Weekly update: insert two new "Betarraga" price rows (Primera / Segunda)
# at the top of the recent-history block (new sheet rows 527-528), pushing
# the existing rows 527-560 down to 529-562.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 527 (shifts 527:560 -> 529:562).
$ws.Range("A527:A528").EntireRow.Insert()

# --- New row 527 : Betarraga, Primera -------------------------------------
$ws.Cells.Item(527, 1).Value  = 1
$ws.Cells.Item(527, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(527, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(527, 4).Value  = 45267
$ws.Cells.Item(527, 5).Value  = 15
$ws.Cells.Item(527, 6).Value  = 100114014
$ws.Cells.Item(527, 7).Value  = "Betarraga"
$ws.Cells.Item(527, 8).Value  = "Sin especificar"
$ws.Cells.Item(527, 9).Value  = "Primera"
$ws.Cells.Item(527, 10).Value = 1200
$ws.Cells.Item(527, 11).Value = 350
$ws.Cells.Item(527, 12).Value = 400
$ws.Cells.Item(527, 13).Value = 375
$ws.Cells.Item(527, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(527, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(527, 16).Value = 94
$ws.Cells.Item(527, 17).Value = 4
$ws.Cells.Item(527, 18).Value = "Hortaliza"

# --- New row 528 : Betarraga, Segunda --------------------------------------
$ws.Cells.Item(528, 1).Value  = 1
$ws.Cells.Item(528, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(528, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(528, 4).Value  = 45267
$ws.Cells.Item(528, 5).Value  = 15
$ws.Cells.Item(528, 6).Value  = 100114014
$ws.Cells.Item(528, 7).Value  = "Betarraga"
$ws.Cells.Item(528, 8).Value  = "Sin especificar"
$ws.Cells.Item(528, 9).Value  = "Segunda"
$ws.Cells.Item(528, 10).Value = 1200
$ws.Cells.Item(528, 11).Value = 350
$ws.Cells.Item(528, 12).Value = 400
$ws.Cells.Item(528, 13).Value = 375
$ws.Cells.Item(528, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(528, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(528, 16).Value = 75
$ws.Cells.Item(528, 17).Value = 5
$ws.Cells.Item(528, 18).Value = "Hortaliza"
